# Add a "Skill Description" column (full skill name) right after the
# existing "SkillCode" column (column A), shifting SFIA Level / Keycode /
# Description one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B (pushes SFIA Level/Keycode/Description right)
$ws.Columns.Item(2).Insert()

# Map each SkillCode (short code) to its full descriptive name
$skillNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "SLEN"       = "Systems and software life cycle engineering"
}

# Header for the new column
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Find the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $skillNames[$code]
}
